$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" and "is_enabled" columns (D:E) from the template,
# shifting the remaining columns (order_by, rem) left.
$ws.Range("D1:E1").EntireColumn.Delete()
